$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 17.65012741088867
$ws.Range("D2").Value = 123

$ws.Range("C3").Value = 17.29297637939453
$ws.Range("D3").Value = 123

$ws.Range("C4").Value = 17.11106300354004
$ws.Range("D4").Value = 123

$ws.Range("C5").Value = 17.25292205810547
$ws.Range("D5").Value = 123

$ws.Range("C6").Value = 17.0590877532959
$ws.Range("D6").Value = 123
